# Apply the "log_prefix_name" row insertion to the Vscs sheet.
# Commit: added schema files from jenkins job
#
# Net effect: a new row is inserted above the current row 29
# ("VSD FQDN ..."), carrying the label "LOG prefix Name" and a comment
# "Log prefix name to pass in the vsc config". Every row from the old
# row 29 through row 99 shifts down by one (to 30..100), along with its
# comment text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vscs")

$firstShiftRow = 29
$lastShiftRow = 99

# 1) Capture the comment text currently attached to column A for every
#    row that is about to shift, keyed by its (pre-insert) row number.
#    Row-insert in this engine moves cell values/styles/validations but
#    leaves comment anchors behind, so we have to move the text by hand.
$commentTexts = @{}
for ($r = $firstShiftRow; $r -le $lastShiftRow; $r++) {
    $cmt = $ws.Range("A" + $r).Comment
    if ($cmt -ne $null) {
        $commentTexts[$r] = $cmt.Text()
    }
}

# 2) Insert a blank row at 29; this pushes rows 29..99 down to 30..100
#    and fixes up dimension/mergeCells/dataValidations/shared formulas
#    automatically.
$ws.Rows.Item($firstShiftRow).Insert()

# 3) The new row 29 inherited style from the row above for columns B/C
#    (plain, no border) instead of the "form field" look used throughout
#    this block. Re-pull the correct formatting from the row directly
#    below (which still carries the original row-29 formatting).
$ws.Range("B" + ($firstShiftRow + 1) + ":C" + ($firstShiftRow + 1)).Copy()
$ws.Range("B" + $firstShiftRow + ":C" + $firstShiftRow).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 4) Remove every old comment still anchored to its pre-insert row
#    (they currently sit one row too high), then re-create them one row
#    lower with the same text, walking bottom-up so we never touch a
#    cell before its own original comment has been captured/cleared.
for ($r = $lastShiftRow; $r -ge $firstShiftRow; $r--) {
    $cmt = $ws.Range("A" + $r).Comment
    if ($cmt -ne $null) {
        $cmt.Delete()
    }
}
for ($r = $lastShiftRow; $r -ge $firstShiftRow; $r--) {
    if ($commentTexts.ContainsKey($r)) {
        $ws.Range("A" + ($r + 1)).AddComment($commentTexts[$r])
    }
}

# 5) Populate the newly-opened row 29 with its label and comment.
$ws.Range("A" + $firstShiftRow).Value = "LOG prefix Name"
$ws.Range("A" + $firstShiftRow).AddComment("Log prefix name to pass in the vsc config")
